# "updated RME jackknives and neg eigen"
#
# 1. Sheet1: the trailing "lamda / dimensionality / lamda dimensionality"
#    summary rows (9-11) are removed (that analysis moved to a new sheet).
# 2. A new "Sheet2" is added holding the RME jackknife ("lamda")
#    neg-eigen ("dimensionality" / "lamda dimensionality") results,
#    split into an "unclustered" block and a "clustered" block, plus an
#    "activity" row for each block.
# 3. Two new (Menlo) fonts are introduced: a light-gray one used for the
#    still-running/placeholder jackknife cells, and a plain black one used
#    for the "activity" rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Sheet1 — drop the old lamda / dimensionality / lamda dimensionality
#    rows (now superseded by Sheet2).
# ---------------------------------------------------------------------
$ws1.Range("A9:I11").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2. Insert the new Sheet2 right after Sheet1.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# 3. Two helper cell styles (Menlo) — built once, applied, then removed
#    again so only the underlying cell formatting (not a named style)
#    survives, matching a hand-formatted-in-Excel workbook.
# ---------------------------------------------------------------------
$grayMenlo = $wb.Styles.Add("GrayMenlo")
$grayMenlo.Font.Name = "Menlo"
$grayMenlo.Font.Color = 13948116   # RGB(212,212,212) -> 0xD4D4D4

$blackMenlo = $wb.Styles.Add("BlackMenlo")
$blackMenlo.Font.Name = "Menlo"

# ---------------------------------------------------------------------
# 4. "unclustered" block (rows 1-6)
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "unclustered"

$hdr = @(0, 1, 2, 3, 4, 5, 6)
for ($i = 0; $i -lt $hdr.Length; $i++) {
    $ws2.Cells.Item(2, 2 + $i).Value = $hdr[$i]
}

$ws2.Range("A3").Value = "lamda"
$ws2.Range("F3").Value = 4.09263
$ws2.Range("G3").Value = 4.11193
$ws2.Range("H3").Value = 3.92698
$ws2.Range("F3:H3").Style = "GrayMenlo"

$ws2.Range("A4").Value = "dimensionality"
$row4 = @(9, 8, 8, 8, 8, 9, 8)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws2.Cells.Item(4, 2 + $i).Value = $row4[$i]
}

$ws2.Range("A5").Value = "lamda dimensionality"
$ws2.Range("E5").Value = 1.9967
$ws2.Range("F5").Value = 1.88625
$ws2.Range("G5").Value = 1.48621
$ws2.Range("H5").Value = 1.32883

$ws2.Range("A6").Value = "activity"
$ws2.Range("B6").Value = 368.87984146231599
$ws2.Range("C6").Value = "408.86759179973296,"
$ws2.Range("D6").Value = 353.98983654932601
$ws2.Range("E6").Value = 389.27098478457702
$ws2.Range("F6").Value = 391.06935931543501
$ws2.Range("G6").Value = 383.10808562598402
$ws2.Range("H6").Value = 399.47930661996099
$ws2.Range("B6").Style = "BlackMenlo"

# ---------------------------------------------------------------------
# 5. "clustered" block (rows 8-12)
# ---------------------------------------------------------------------
$ws2.Range("A8").Value = "clustered"
for ($i = 0; $i -lt $hdr.Length; $i++) {
    $ws2.Cells.Item(8, 2 + $i).Value = $hdr[$i]
}

$ws2.Range("A9").Value = "lamda"
$ws2.Range("B9").Value = 3.59897
$ws2.Range("C9").Value = 3.00975
$ws2.Range("D9").Value = 2.89703
$ws2.Range("E9").Value = 3.14789
$ws2.Range("F9").Value = 2.99823
$ws2.Range("G9").Value = 3.0982
$ws2.Range("H9").Value = 3.31061
$ws2.Range("B9").Style = "BlackMenlo"

$ws2.Range("A10").Value = "dimensionality"
$row10 = @(4, 4, 4, 3, 4, 4, 4)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws2.Cells.Item(10, 2 + $i).Value = $row10[$i]
}

$ws2.Range("A11").Value = "lamda dimensionality"
$ws2.Range("B11").Value = 2.64442
$ws2.Range("C11").Value = 2.44122
$ws2.Range("D11").Value = 2.11953
$ws2.Range("E11").Value = 3.2202
$ws2.Range("F11").Value = 2.03097
$ws2.Range("G11").Value = 2.44319
$ws2.Range("H11").Value = 2.59429

$ws2.Range("A12").Value = "activity"

# Remove the helper named styles again — the formatting they carried
# stays on the cells, only the style *definitions* go away.
$grayMenlo.Delete()
$blackMenlo.Delete()

# ---------------------------------------------------------------------
# 6. Column width + selections + which sheet/tab is active.
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 22.1640625

$ws1.Range("A9:B11").Select()
$ws2.Range("H11").Select()
$ws2.Activate()
